# Rename the existing sheet and add a new "EventsChart" sheet after it,
# then populate it with the Oregon Trail "events" table.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "RidersCalc"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "EventsChart"

# Header row
$ws2.Range("A1").Value = "Odds"
$ws2.Range("B1").Value = "Event #"
$ws2.Range("C1").Value = "Event Line #"
$ws2.Range("D1").WrapText = $true

# Data rows: Odds, Event #, Event Line #
$data = @(
    @(6, 1, 3660),
    @(11, 2, 3700),
    @(13, 3, 3740),
    @(15, 4, 3790),
    @(17, 5, 3820),
    @(22, 6, 3850),
    @(32, 7, 3880),
    @(35, 8, 3960),
    @(37, 9, 4130),
    @(42, 10, 4190),
    @(44, 11, 4220),
    @(54, 12, 4290),
    @(64, 13, 4340),
    @(69, 14, 4650),
    @(95, 15, 4610)
)

$row = 2
foreach ($item in $data) {
    $ws2.Cells.Item($row, 1).Value = $item[0]
    $ws2.Cells.Item($row, 2).Value = $item[1]
    $ws2.Cells.Item($row, 3).Value = $item[2]
    $row = $row + 1
}

# Final row has no "Odds" value, only Event # and Event Line #
$ws2.Cells.Item($row, 2).Value = 16
$ws2.Cells.Item($row, 3).Value = 4670

# Match the recorded selection on the new sheet
$ws2.Range("B18").Select() | Out-Null
